$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the test data values (shared strings POC5OCTo -> POC15OCTo, etc.)
$ws.Range("A5").Value = "POC15OCTo"
$ws.Range("B5").Value = "POC15OCTo"
$ws.Range("C5").Value = "FacilityPOC15OCTo"
$ws.Range("D5").Value = "FacilityPOC15OCTo"
$ws.Range("E5").Value = "PharmacyPOC15OCTo"
$ws.Range("F5").Value = "PharmacyPOC15OCTo"

# Update the active selection to match the new cursor position
$ws.Range("H9").Select()
